$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.648.06"
$ws.Range("E2").Value = "  +3.91%  "

$ws.Range("D3").Value = "2.631.81"
$ws.Range("E3").Value = "  +4.39%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'606.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.05%  "

$ws.Range("D6").Value = "'179.26"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.55%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +1.56%  "

$ws.Range("D9").Value = "2.631.93"
$ws.Range("E9").Value = "  +4.33%  "

$ws.Range("D10").Value = "'0.169"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +14.78%  "

$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("E12").Value = "  +2.85%  "

$ws.Range("D13").Value = "'5.05"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("D14").Value = "3.136.08"
$ws.Range("E14").Value = "  +4.26%  "

$ws.Range("D15").Value = "'0.0000187"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +8.48%  "

$ws.Range("D16").Value = "'26.59"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.64%  "

$ws.Range("D17").Value = "71.601.05"
$ws.Range("E17").Value = "  +4.00%  "

$ws.Range("D18").Value = "2.640.83"
$ws.Range("E18").Value = "  +5.57%  "

$ws.Range("D19").Value = "'383.00"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.64%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'7.98"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.17%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'11.46"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.02%  "

$ws.Range("D22").Value = "'4.13"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.22%  "

$ws.Range("D23").Value = "'1.98"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +17.21%  "

$ws.Range("B24").Value = "NEARProtocol"
$ws.Range("C24").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D24").Value = "'4.50"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.21%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'72.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.73%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").Value = "'10.02"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +10.94%  "

$ws.Range("D28").Value = "2.766.36"
$ws.Range("E28").Value = "  +4.44%  "

$ws.Range("D29").Value = "'0.996"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("D30").Value = "'550.43"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.48%  "

$ws.Range("E31").Value = "  +8.78%  "

$ws.Range("D32").Value = "'8.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.93%  "

$ws.Range("E33").Value = "  +7.64%  "

$ws.Range("E34").Value = "  +3.13%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").Value = "'166.45"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.93%  "

$ws.Range("E37").Value = "  +3.11%  "

$ws.Range("D38").Value = "'0.114"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.73%  "

$ws.Range("D39").Value = "'19.18"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.74%  "

$ws.Range("D40").Value = "'1.40"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.63%  "

$ws.Range("D41").Value = "'1.87"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.29%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.63"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +11.05%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").Value = "'5.04"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.90%  "

$ws.Range("D45").Value = "'0.333"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.22%  "

$ws.Range("D46").Value = "'39.19"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("D47").Value = "'152.01"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.58%  "

$ws.Range("D48").Value = "'3.64"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.83%  "

$ws.Range("D49").Value = "'0.536"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.79%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'6.29"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +10.79%  "

$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").Value = "'1.69"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.27%  "
